$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Delete the "Desarquivamentos Pendentes" sheet entirely
$ws = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$ws.Delete()

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$ws1 = $wb.Worksheets.Item("Paineis DARQ")
$ws1.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$ws2 = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$ws2.Name = "RECOLHIMENTO X ELIMINAÇÃO"
